# Log Week 17 data across the Football Team stat sheets and fix the
# Simulate_Season.py tiebreaking totals that depend on it.
#
# Each "log" cell on the YDS / ST sheets stores one number per game played
# so far as a space-separated string; logging a new week just appends the
# new game's number(s) to the end of the existing string. The season-total
# cells on OFF / DEF / ST / TURNS / PEN simply get bumped by the Week 17
# numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append Week 17 per-game yardage logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$r = $ws.Range("B2")
$r.Value = $r.Value() + " 3 3 11 3 6 1 1 3 4 4 5 8 1 4 3 8 0 4 10 1 8 3"

$r = $ws.Range("C2")
$r.Value = $r.Value() + " 0 6 6 -3 7 4 3 0 3 -3 22 3 6 4 2 0 1 14 13 5 4 2 8 0 1 7 4 2 3"

$r = $ws.Range("B3")
$r.Value = $r.Value() + " 16 10 24 8 25 12 2 12 12 11 5 6 18 9 4 11 1 -4 2 8 10 4 5 9 12 5 10"

$r = $ws.Range("C3")
$r.Value = $r.Value() + " 30 5 12 14 5 19 18 5 8 8 18 13 -3 27 12 13 10"

# ---------------------------------------------------------------------
# OFF sheet: season totals through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = 190
$ws.Range("D2").Value = 10
$ws.Range("F2").Value = 77
$ws.Range("G2").Value = 53
$ws.Range("J2").Value = 44
$ws.Range("L2").Value = 287
$ws.Range("M2").Value = 190
$ws.Range("O2").Value = 31
$ws.Range("Q2").Value = 534

$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 196
$ws.Range("E3").Value = 33
$ws.Range("F3").Value = 92
$ws.Range("H3").Value = 28
$ws.Range("I3").Value = 58
$ws.Range("J3").Value = 54
$ws.Range("N3").Value = 18

# ---------------------------------------------------------------------
# DEF sheet: season totals through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("C2").Value = 181
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 65
$ws.Range("G2").Value = 44
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 33
$ws.Range("L2").Value = 287
$ws.Range("M2").Value = 183
$ws.Range("O2").Value = 18
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 503

$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 194
$ws.Range("F3").Value = 120
$ws.Range("G3").Value = 38
$ws.Range("H3").Value = 22
$ws.Range("I3").Value = 52
$ws.Range("J3").Value = 72
$ws.Range("N3").Value = 18

# ---------------------------------------------------------------------
# ST sheet: season totals + per-game logs through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 72
$ws.Range("D2").Value = 56
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 87
$ws.Range("J2").Value = 47
$ws.Range("K2").Value = 45
$ws.Range("N2").Value = 21
$ws.Range("O2").Value = 13

$ws.Range("B3").Value = 46

$r = $ws.Range("B4")
$r.Value = $r.Value() + " 66 38"

$r = $ws.Range("B5")
$r.Value = $r.Value() + " 21 0"

$r = $ws.Range("B6")
$r.Value = $r.Value() + " 34 23"

$r = $ws.Range("D3")
$r.Value = $r.Value() + " 38 46"

$r = $ws.Range("D4")
$r.Value = $r.Value() + " 0 7"

$r = $ws.Range("D5")
$r.Value = $r.Value() + " 0 0"

# ---------------------------------------------------------------------
# TURNS sheet: season totals through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("B2").Value = 7
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 6

# ---------------------------------------------------------------------
# PEN sheet: season totals through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B2").Value = 15
